$d = $word.ActiveDocument

# --- Change 1 ---------------------------------------------------------
# "# Data Path: Entities Folder -> 1-Entities Database -> 1-Persons
# Table -> ID: ..." is boilerplate that repeats several times through
# the document, but only the very first occurrence (the one wrapped by
# bookmark "_Hlk210618795", right under the doc title) gets the leading
# "1-" removed from "Entities Database" / "Persons Table". Scope the
# Find/Replace to that bookmark's range so the other repeated copies
# are left untouched.
$bm = $d.Bookmarks("_Hlk210618795")
$r1 = $bm.Range
$r1.Find.ClearFormatting()
$r1.Find.Replacement.ClearFormatting()
$r1.Find.Execute(
    "1-Entities Database -> 1-Persons Table",  # FindText
    $true,                                      # MatchCase
    $false,                                     # MatchWholeWord
    $false,                                     # MatchWildcards
    $false,                                     # MatchSoundsLike
    $false,                                     # MatchAllWordForms
    $true,                                      # Forward
    1,                                           # Wrap (wdFindContinue)
    $false,                                     # Format
    "Entities Database -> Persons Table",        # ReplaceWith
    2                                            # Replace (wdReplaceOne)
)

# --- Change 2 ---------------------------------------------------------
# "    Include bug type and UI location in the TaskName. Use the
# format:" was stored as three separate runs, with a pair of
# <w:proofErr> spell-check markers bracketing the "TaskName" run.
# Re-typing the sentence as one contiguous string collapses it back
# into a single run and drops the now-superfluous proofErr markers.
# This phrase only occurs once in the document, so it is safe to
# search the whole story.
$r2 = $d.Content
$r2.Find.ClearFormatting()
$r2.Find.Replacement.ClearFormatting()
$r2.Find.Execute(
    "Include bug type and UI location in the TaskName. Use the format:",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Include bug type and UI location in the TaskName. Use the format:",
    2
)
